# Leave Card update — 12/22/2023 10:59 AM
# Applies the SICK LEAVE table additions (rows 365-371) and the resulting
# one-row insertion in column A (PERIOD) of Table1, plus the new trailing
# row 512 that the column-A shift pushes out of the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# ---------------------------------------------------------------------
# 1) New SICK LEAVE entries, rows 365-369: EARNED (C) = 1.25 each.
#    The table's calculated column (G) re-derives itself automatically.
# ---------------------------------------------------------------------
foreach ($r in 365..369) {
    $ws.Cells.Item($r, 3).Value = 1.25   # column C = EARNED
}

# Row 367 sits on a heavier "quarter" border (style 41/42/15/12) — typing
# into it the way the original author did (paste of an adjoining cell)
# carried the neighbouring row's border-less number style into column C
# only, while the calculated column G keeps its original bordered style.
$ws.Range("C366").Copy() | Out-Null
$ws.Range("C367").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(367, 3).Value = 1.25
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Row 370 becomes a "SL(7-0-0)" entry: EARNED=1.25, Absence W/Pay=7,
#    REMARKS = "9/30 - 10/10/2023".
# ---------------------------------------------------------------------
$ws.Cells.Item(370, 2).Value = "SL(7-0-0)"            # B370 PARTICULARS
$ws.Cells.Item(370, 3).Value = 1.25                    # C370 EARNED
$ws.Cells.Item(370, 8).Value = 7                       # H370 Absence W/Pay
$ws.Cells.Item(370, 11).Value = "9/30 - 10/10/2023"    # K370 REMARKS

# ---------------------------------------------------------------------
# 3) Row 371 becomes a "SL(4-0-0)" entry: Absence W/Pay=4,
#    REMARKS = "10/17-20/2023". Its PERIOD (A371) date is vacated — the
#    date sequence shifts down by one row from here on (see step 4).
# ---------------------------------------------------------------------
$ws.Cells.Item(371, 2).Value = "SL(4-0-0)"             # B371 PARTICULARS
$ws.Cells.Item(371, 8).Value = 4                       # H371 Absence W/Pay
$ws.Cells.Item(371, 11).Value = "10/17-20/2023"        # K371 REMARKS

# ---------------------------------------------------------------------
# 4) Append one row to the table so the PERIOD column has somewhere to
#    shift into, then push every PERIOD date from A371 down by one row
#    (A371 itself goes blank), working from the bottom up so no value
#    is overwritten before it is copied forward.
# ---------------------------------------------------------------------
$lo.ListRows.Add() | Out-Null
$ws.Range("A511:K511").Copy($ws.Range("A512")) | Out-Null
$excel.CutCopyMode = 0
$ws.Cells.Item(512, 7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

for ($r = 511; $r -ge 372; $r--) {
    $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($r - 1, 1).Value2
}
$ws.Cells.Item(371, 1).Value2 = ""

# ---------------------------------------------------------------------
# 5) Park the cursor where the author left it.
# ---------------------------------------------------------------------
$ws.Range("K371").Select() | Out-Null
